$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = "128HD24COMBO"
$ws.Range("B6").Value = "Container - HD Deli (24oz) w/ Lid"
$ws.Range("C6").Value = "'2"
$ws.Range("D6").Value = "'39.99"
$ws.Range("E6").Value = "'79.98"

# Row 7
$ws.Range("A7").Value = "128HD8BULK"
$ws.Range("B7").Value = "Container - HD Deli (8oz)"
$ws.Range("C7").Value = "'2"
$ws.Range("D7").Value = "'26.49"
$ws.Range("E7").Value = "'52.98"
